# Updated cryptos list on Thu Feb  1 07:48:33 UTC 2024 with GitHub Actions
#
# Refreshes the "Price" (D) and "Volume(1h)" (E) columns for every coin row
# (rows 2-51) with newly scraped values. Price values that look like plain
# numbers (e.g. "95.44") are written with a leading apostrophe so Excel
# keeps them as literal text (matching the source data's exact formatting,
# including trailing zeros such as "34.00" or "0.0990") instead of coercing
# them into numeric cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Text)
    # A leading apostrophe is Excel's standard "treat as text" quote-prefix.
    # It is not stored as part of the cell's value.
    $Range.Value = "'" + $Text
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "42.269.70"
$ws.Range("E2").Value = "  -1.75%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.273.44"
$ws.Range("E3").Value = "  -2.79%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 - BNB
$ws.Range("E5").Value = "  -2.77%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "95.44"
$ws.Range("E6").Value = "  -5.53%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.03%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -3.78%  "

# Row 9 - Cardano
Set-TextValue $ws.Range("D9") "0.492"
$ws.Range("E9").Value = "  -3.63%  "

# Row 10 - Avalanche
Set-TextValue $ws.Range("D10") "33.39"
$ws.Range("E10").Value = "  -4.33%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -0.99%  "

# Row 12 - OKB
Set-TextValue $ws.Range("D12") "48.48"
$ws.Range("E12").Value = "  -7.62%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +0.02%  "

# Row 14 - Polkadot
Set-TextValue $ws.Range("D14") "6.66"
$ws.Range("E14").Value = "  -3.06%  "

# Row 15 - Chainlink
Set-TextValue $ws.Range("D15") "15.70"
$ws.Range("E15").Value = "  -0.63%  "

# Row 16 - Wrapped liquid staked Ether 2.0
$ws.Range("D16").Value = "2.627.44"
$ws.Range("E16").Value = "  -2.91%  "

# Row 17 - Wrapped Ether
$ws.Range("D17").Value = "2.282.23"
$ws.Range("E17").Value = "  -3.37%  "

# Row 18 - Polygon
$ws.Range("E18").Value = "  -5.78%  "

# Row 19 - (42.xxx.xx priced coin)
$ws.Range("D19").Value = "42.211.92"
$ws.Range("E19").Value = "  -1.72%  "

# Row 20
Set-TextValue $ws.Range("D20") "11.63"
$ws.Range("E20").Value = "  -0.87%  "

# Row 21
$ws.Range("E21").Value = "  -2.02%  "

# Row 22
$ws.Range("E22").Value = "  -3.77%  "

# Row 23
Set-TextValue $ws.Range("D23") "66.72"
$ws.Range("E23").Value = "  -1.93%  "

# Row 24
Set-TextValue $ws.Range("D24") "233.51"
$ws.Range("E24").Value = "  -1.31%  "

# Row 25
$ws.Range("E25").Value = "  -2.43%  "

# Row 26
$ws.Range("E26").Value = "  +0.24%  "

# Row 27
Set-TextValue $ws.Range("D27") "2.46"
$ws.Range("E27").Value = "  -4.09%  "

# Row 28
Set-TextValue $ws.Range("D28") "23.90"
$ws.Range("E28").Value = "  -6.43%  "

# Row 29
$ws.Range("E29").Value = "  -1.06%  "

# Row 30
Set-TextValue $ws.Range("D30") "167.61"
$ws.Range("E30").Value = "  +2.24%  "

# Row 31
Set-TextValue $ws.Range("D31") "34.00"
$ws.Range("E31").Value = "  -4.79%  "

# Row 32
Set-TextValue $ws.Range("D32") "9.10"
$ws.Range("E32").Value = "  -2.20%  "

# Row 33
$ws.Range("E33").Value = "  -0.05%  "

# Row 34
$ws.Range("E34").Value = "  -3.94%  "

# Row 35
$ws.Range("E35").Value = "  -2.58%  "

# Row 36
$ws.Range("E36").Value = "  -4.85%  "

# Row 37
$ws.Range("E37").Value = "  -4.87%  "

# Row 38
Set-TextValue $ws.Range("D38") "16.39"
$ws.Range("E38").Value = "  -6.58%  "

# Row 39
$ws.Range("E39").Value = "  -4.17%  "

# Row 40
Set-TextValue $ws.Range("D40") "0.0990"
$ws.Range("E40").Value = "  -2.99%  "

# Row 41
$ws.Range("E41").Value = "  -3.16%  "

# Row 42
$ws.Range("E42").Value = "  -6.78%  "

# Row 43
$ws.Range("E43").Value = "  -6.81%  "

# Row 44 - Maker
$ws.Range("D44").Value = "1.965.20"
$ws.Range("E44").Value = "  -2.99%  "

# Row 45
$ws.Range("E45").Value = "  -2.18%  "

# Row 46
Set-TextValue $ws.Range("D46") "17.57"
$ws.Range("E46").Value = "  -7.20%  "

# Row 47
Set-TextValue $ws.Range("D47") "9.61"
$ws.Range("E47").Value = "  -5.58%  "

# Row 48
$ws.Range("E48").Value = "  -4.73%  "

# Row 49 - RocketPoolETH
$ws.Range("D49").Value = "2.498.83"
$ws.Range("E49").Value = "  -2.33%  "

# Row 50
Set-TextValue $ws.Range("D50") "52.33"
$ws.Range("E50").Value = "  -7.70%  "

# Row 51
$ws.Range("E51").Value = "  -5.58%  "
